# Updates the cryptos list (Price / Volume(1h) columns, and the
# Frax/EnergySwap row swap) to match the latest scrape, per the
# "Updated cryptos list ... with GitHub Actions" commit.
#
# Note: a few Price values (e.g. "6.220", "9.980", "125.60",
# "0.00001134") look numeric but must stay literal text (trailing
# zeros / exact digit count matter), so they're entered with a
# leading single-quote to force Excel to store them as text instead
# of auto-converting to a Double and losing the formatting.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "22.425.10"
$ws.Range("E2").Value = "  +0.10%  "
$ws.Range("D3").Value = "1.572.30"
$ws.Range("E3").Value = "  +0.12%  "
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("E5").Value = "  +0.18%  "
$ws.Range("D6").Value = "291.16"
$ws.Range("E6").Value = "  +0.05%  "
$ws.Range("D7").Value = "0.3764"
$ws.Range("E7").Value = "  +2.08%  "
$ws.Range("D8").Value = "49.92"
$ws.Range("E8").Value = "  +0.59%  "
$ws.Range("D9").Value = "0.3415"
$ws.Range("E9").Value = "  +1.48%  "
$ws.Range("D10").Value = "1.162"
$ws.Range("E10").Value = "  -0.11%  "
$ws.Range("D11").Value = "0.07665"
$ws.Range("E11").Value = "  +1.37%  "
$ws.Range("D12").Value = "1.002"
$ws.Range("E12").Value = "  +0.13%  "
$ws.Range("D13").Value = "21.25"
$ws.Range("E13").Value = "  +1.02%  "
$ws.Range("D14").Value = "5.976"
$ws.Range("E14").Value = "  -1.28%  "
$ws.Range("D15").Value = "6.915"
$ws.Range("E15").Value = "  +1.01%  "
$ws.Range("D16").Value = "1.572.37"
$ws.Range("E16").Value = "  -0.10%  "
$ws.Range("D17").Value = "'0.00001134"
$ws.Range("E17").Value = "  -0.50%  "
$ws.Range("D18").Value = "90.29"
$ws.Range("E18").Value = "  +1.21%  "
$ws.Range("D19").Value = "0.06775"
$ws.Range("E19").Value = "  +1.39%  "
$ws.Range("E20").Value = "  +0.28%  "
$ws.Range("D21").Value = "16.79"
$ws.Range("E21").Value = "  +2.89%  "
$ws.Range("D22").Value = "'6.220"
$ws.Range("E22").Value = "  -0.20%  "
$ws.Range("D23").Value = "0.5271"
$ws.Range("E23").Value = "  -2.79%  "
$ws.Range("D24").Value = "12.04"
$ws.Range("E24").Value = "  +0.76%  "
$ws.Range("D25").Value = "2.422"
$ws.Range("E25").Value = "  +0.89%  "
$ws.Range("D26").Value = "22.402.33"
$ws.Range("E26").Value = "  -0.05%  "
$ws.Range("D27").Value = "2.729"
$ws.Range("E27").Value = "  -7.29%  "
$ws.Range("E28").Value = "  +2.28%  "
$ws.Range("D29").Value = "145.17"
$ws.Range("E29").Value = "  -0.72%  "
$ws.Range("D30").Value = "5.061"
$ws.Range("E30").Value = "  +2.65%  "
$ws.Range("D31").Value = "125.97"
$ws.Range("E31").Value = "  +0.74%  "
$ws.Range("D32").Value = "1.744.28"
$ws.Range("E32").Value = "  -0.30%  "
$ws.Range("D33").Value = "'6.200"
$ws.Range("E33").Value = "  -0.90%  "
$ws.Range("D34").Value = "1.009"
$ws.Range("E34").Value = "  +3.18%  "
$ws.Range("D35").Value = "2.009"
$ws.Range("E35").Value = "  +1.99%  "
$ws.Range("D36").Value = "'9.980"
$ws.Range("E36").Value = "  -3.48%  "
$ws.Range("D37").Value = "0.08555"
$ws.Range("E37").Value = "  +0.74%  "
$ws.Range("D38").Value = "0.02551"
$ws.Range("E38").Value = "  +0.92%  "
$ws.Range("E39").Value = "  +0.76%  "
$ws.Range("D40").Value = "1.347"
$ws.Range("E40").Value = "  +7.91%  "
$ws.Range("D41").Value = "0.06536"
$ws.Range("E41").Value = "  +0.52%  "
$ws.Range("D42").Value = "5.461"
$ws.Range("E42").Value = "  -0.32%  "
$ws.Range("D43").Value = "0.6457"
$ws.Range("E43").Value = "  +1.14%  "
$ws.Range("D44").Value = "11.58"
$ws.Range("E44").Value = "  -1.60%  "
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").Value = "14.13"
$ws.Range("E45").Value = "  -2.58%  "
$ws.Range("B46").Value = "Frax"
$ws.Range("C46").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D46").Value = "1.002"
$ws.Range("E46").Value = "  +0.20%  "
$ws.Range("D47").Value = "0.6027"
$ws.Range("E47").Value = "  +0.65%  "
$ws.Range("D48").Value = "3.788"
$ws.Range("E48").Value = "  +0.51%  "
$ws.Range("D49").Value = "1.298"
$ws.Range("E49").Value = "  +9.07%  "
$ws.Range("D50").Value = "2.091"
$ws.Range("E50").Value = "  -0.45%  "
$ws.Range("D51").Value = "'125.60"
$ws.Range("E51").Value = "  +3.13%  "
